$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "On Order" column (I) ---------------------------------------------
# Header
$ws.Range("I1").Value = "On Order"

# Mark every existing BOM line "Yes" in the new column, mirroring the
# existing "Footprint Checked" (H) column that already has "Yes" for every
# row except rows 11-14 (which instead hold N/A in H, but still get a Yes
# in the new On Order column per the diff).
foreach ($r in 2..36) {
    $ws.Cells.Item($r, 9).Value = "Yes"
}

# --- New BOM line: row 36 now also gets a part number / manufacturer -------
$ws.Range("A36").Value = "RC1005J000CS"
$ws.Range("B36").Value = "Samsung"

# Give A36 the same "quoted hyperlink-style" look used by the other part
# links in column A (style index 8 in the original workbook).
$ws.Range("A36").Style = $ws.Range("A18").Style

# Add the Digikey hyperlink for the new part, matching the style used by
# the other hyperlinks in column A.
$ws.Hyperlinks.Add($ws.Range("A36"), "https://www.digikey.com/product-detail/en/samsung-electro-mechanics/RC1005J000CS/1276-3480-1-ND/3903583", "", "", "https://www.digikey.com/product-detail/en/samsung-electro-mechanics/RC1005J000CS/1276-3480-1-ND/3903583")

# --- View state: scroll down a bit and move the selection ------------------
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("K7").Select()
